# re-run RU 1001; without crop
# Set Saudi Arabia column (C) values for rows 2-6 to #NUM! errors,
# matching a re-run of the underlying data source that failed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "#NUM!"
$ws.Range("C3").Value = "#NUM!"
$ws.Range("C4").Value = "#NUM!"
$ws.Range("C5").Value = "#NUM!"
$ws.Range("C6").Value = "#NUM!"
